$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-29 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-30 Saturday", 2) | Out-Null
$d.Content.Find.Execute("419÷2=209, 1", $true, $false, $false, $false, $false, $true, 1, $false, "999÷4=249, 3", 2) | Out-Null
$d.Content.Find.Execute("574÷5=114, 4", $true, $false, $false, $false, $false, $true, 1, $false, "829÷5=165, 4", 2) | Out-Null
$d.Content.Find.Execute("783÷9=87, 0", $true, $false, $false, $false, $false, $true, 1, $false, "793÷9=88, 1", 2) | Out-Null
$d.Content.Find.Execute("878÷6=146, 2", $true, $false, $false, $false, $false, $true, 1, $false, "620÷7=88, 4", 2) | Out-Null
$d.Content.Find.Execute("637÷2=318, 1", $true, $false, $false, $false, $false, $true, 1, $false, "249÷6=41, 3", 2) | Out-Null
$d.Content.Find.Execute("677÷4=169, 1", $true, $false, $false, $false, $false, $true, 1, $false, "132÷5=26, 2", 2) | Out-Null
$d.Content.Find.Execute("613÷3=204, 1", $true, $false, $false, $false, $false, $true, 1, $false, "623÷8=77, 7", 2) | Out-Null
$d.Content.Find.Execute("500÷6=83, 2", $true, $false, $false, $false, $false, $true, 1, $false, "528÷9=58, 6", 2) | Out-Null
$d.Content.Find.Execute("167÷9=18, 5", $true, $false, $false, $false, $false, $true, 1, $false, "465÷5=93, 0", 2) | Out-Null
$d.Content.Find.Execute("337÷9=37, 4", $true, $false, $false, $false, $false, $true, 1, $false, "335÷8=41, 7", 2) | Out-Null
$d.Content.Find.Execute("359÷3=119, 2", $true, $false, $false, $false, $false, $true, 1, $false, "313÷6=52, 1", 2) | Out-Null
$d.Content.Find.Execute("503÷5=100, 3", $true, $false, $false, $false, $false, $true, 1, $false, "956÷7=136, 4", 2) | Out-Null
$d.Content.Find.Execute("442÷7=63, 1", $true, $false, $false, $false, $false, $true, 1, $false, "442÷4=110, 2", 2) | Out-Null
$d.Content.Find.Execute("252÷7=36, 0", $true, $false, $false, $false, $false, $true, 1, $false, "227÷9=25, 2", 2) | Out-Null
$d.Content.Find.Execute("784÷3=261, 1", $true, $false, $false, $false, $false, $true, 1, $false, "471÷6=78, 3", 2) | Out-Null
$d.Content.Find.Execute("192÷6=32, 0", $true, $false, $false, $false, $false, $true, 1, $false, "267÷9=29, 6", 2) | Out-Null
$d.Content.Find.Execute("117÷9=13, 0", $true, $false, $false, $false, $false, $true, 1, $false, "858÷9=95, 3", 2) | Out-Null
$d.Content.Find.Execute("702÷9=78, 0", $true, $false, $false, $false, $false, $true, 1, $false, "583÷3=194, 1", 2) | Out-Null
$d.Content.Find.Execute("257÷2=128, 1", $true, $false, $false, $false, $false, $true, 1, $false, "564÷2=282, 0", 2) | Out-Null
$d.Content.Find.Execute("856÷2=428, 0", $true, $false, $false, $false, $false, $true, 1, $false, "824÷6=137, 2", 2) | Out-Null
$d.Content.Find.Execute("654÷7=93, 3", $true, $false, $false, $false, $false, $true, 1, $false, "390÷6=65, 0", 2) | Out-Null
$d.Content.Find.Execute("266÷8=33, 2", $true, $false, $false, $false, $false, $true, 1, $false, "290÷9=32, 2", 2) | Out-Null
$d.Content.Find.Execute("480÷3=160, 0", $true, $false, $false, $false, $false, $true, 1, $false, "580÷2=290, 0", 2) | Out-Null
$d.Content.Find.Execute("320÷6=53, 2", $true, $false, $false, $false, $false, $true, 1, $false, "556÷3=185, 1", 2) | Out-Null
$d.Content.Find.Execute("257÷8=32, 1", $true, $false, $false, $false, $false, $true, 1, $false, "303÷4=75, 3", 2) | Out-Null
